$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/coverage-days"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Clear the root Extension row's Constraint(s) value (column AI, row 2)
$elem.Range("AI2").Value = ""

# The Extension.url row's Fixed Value mirrors the StructureDefinition's own
# canonical URL (shared string also used by Metadata!B2), so it must track
# the new URL too.
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/coverage-days"
